$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet (sheet1): insert a new 2022-Q3 summary row above the
#    existing 2022-Q2 row (so 2022-Q3 becomes row 2, 2022-Q2 shifts to row 3).
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Clone A2's cell format onto A3 first (brand-new row) so the shifted
# 2022-Q2 row keeps the same "index" style as the header/label column.
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

# Move the current 2022-Q2 figures down into row 3.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 1
$wsTotal.Range("D3").Value = 0

# Overwrite row 2 with the new 2022-Q3 figures.
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.19

# ---------------------------------------------------------------------------
# 2) Duplicate the existing "2022-Q2" sheet so its current (Q2) data is
#    preserved on a new tab, then rename the original tab to "2022-Q3" and
#    overwrite its contents with the fresh Q3 fund data.
# ---------------------------------------------------------------------------
$wsQ2Old = $wb.Worksheets.Item(2)

# Copy the sheet to just after itself -> becomes the new "2022-Q2" tab.
$wsQ2Old.Copy($null, $wsQ2Old)

$wsQ3 = $wsQ2Old
$wsQ3.Name = "2022-Q3"

$wsQ2New = $wb.Worksheets.Item(3)
$wsQ2New.Name = "2022-Q2"

# ---------------------------------------------------------------------------
# 3) Re-style the "2022-Q3" sheet's header row and index column to style 2
#    (matching the rest of the workbook) by copying formatting from the
#    "总计" sheet, which already uses that style - this re-uses the existing
#    style table entry instead of creating a new one.
# ---------------------------------------------------------------------------
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Write the 2022-Q3 fund data. Numeric-looking text (fund codes with
#    leading zeros, percentage/size figures stored as text in the source
#    data) is entered with a leading apostrophe so it is kept as text
#    instead of being auto-converted to a number.
# ---------------------------------------------------------------------------
$wsQ3.Range("A2").Value = 0
$wsQ3.Range("B2").Value = "'005433"
$wsQ3.Range("C2").Value = "申万菱信医药先锋股票A"
$wsQ3.Range("D2").Value = "'1.76"
$wsQ3.Range("E2").Value = "'93.02"
$wsQ3.Range("F2").Value = "'8.07"
$wsQ3.Range("G2").Value = "'0.1420"
$wsQ3.Range("H2").Value = 2

$wsQ3.Range("A3").Value = 1
$wsQ3.Range("B3").Value = "'014220"
$wsQ3.Range("C3").Value = "恒越医疗健康精选混合A"
$wsQ3.Range("D3").Value = "'0.72"
$wsQ3.Range("E3").Value = "'88.76"
$wsQ3.Range("F3").Value = "'4.73"
$wsQ3.Range("G3").Value = "'0.0341"
$wsQ3.Range("H3").Value = 6

$wsQ3.Range("A4").Value = 2
$wsQ3.Range("B4").Value = "'014221"
$wsQ3.Range("C4").Value = "恒越医疗健康精选混合C"
$wsQ3.Range("D4").Value = "'0.29"
$wsQ3.Range("E4").Value = "'88.76"
$wsQ3.Range("F4").Value = "'4.73"
$wsQ3.Range("G4").Value = "'0.0137"
$wsQ3.Range("H4").Value = 6

$wsQ3.Range("A5").Value = 3
$wsQ3.Range("B5").Value = "'015171"
$wsQ3.Range("C5").Value = "申万菱信医药先锋股票C"
$wsQ3.Range("D5").Value = "'0.00"
$wsQ3.Range("E5").Value = "'93.02"
$wsQ3.Range("F5").Value = "'8.07"
$wsQ3.Range("G5").Value = 0
$wsQ3.Range("H5").Value = 2
